$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 770
$ws.Range("F6").Value = 2494
$ws.Range("F8").Value = 1842
$ws.Range("F9").Value = 3166
$ws.Range("F11").Value = 4657
$ws.Range("F12").Value = 432
$ws.Range("F13").Value = 254
$ws.Range("F14").Value = 148
$ws.Range("F15").Value = 598
$ws.Range("F17").Value = 5
$ws.Range("F20").Value = 274
$ws.Range("F21").Value = 13
$ws.Range("F23").Value = 130
$ws.Range("F24").Value = 325
$ws.Range("F25").Value = 4644
$ws.Range("F26").Value = 10
$ws.Range("F27").Value = 27
$ws.Range("F29").Value = 5188
$ws.Range("F31").Value = 1165
$ws.Range("F32").Value = 227
$ws.Range("F33").Value = 635
$ws.Range("F35").Value = 1
$ws.Range("F36").Value = 61
$ws.Range("F38").Value = 754
$ws.Range("F39").Value = 48
$ws.Range("F40").Value = 690
$ws.Range("F41").Value = 690

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 1064
$ws.Range("F4").Value = 28

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 1064
$ws.Range("F5").Value = 28
$ws.Range("F8").Value = 770
$ws.Range("F9").Value = 2494
$ws.Range("F11").Value = 1842
$ws.Range("F13").Value = 3166
$ws.Range("F15").Value = 4657
$ws.Range("F16").Value = 432
$ws.Range("F17").Value = 254
$ws.Range("F18").Value = 148
$ws.Range("F19").Value = 598
$ws.Range("F21").Value = 5
$ws.Range("F24").Value = 274
$ws.Range("F25").Value = 13
$ws.Range("F28").Value = 130
$ws.Range("F29").Value = 324
$ws.Range("F30").Value = 4644
$ws.Range("F31").Value = 10
$ws.Range("F32").Value = 27
$ws.Range("F34").Value = 5188
$ws.Range("F36").Value = 1165
$ws.Range("F37").Value = 227
$ws.Range("F38").Value = 635
$ws.Range("F40").Value = 1
$ws.Range("F41").Value = 7
$ws.Range("F42").Value = 61
$ws.Range("F44").Value = 754
$ws.Range("F45").Value = 48
$ws.Range("F46").Value = 689
$ws.Range("F47").Value = 690
